# Daily attendance processing - 2026-01-31 15:38:18
# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (column G) wherever they appear together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$col = $ws.Range("G1:G$lastRow")

$col.Replace($oldValue, $newValue)
